# Apply the "Add files via upload" revision to the Veeva manual-testing
# workbook: fill in the previously-empty "Test Scenario Description"
# column (B) for the generic testing-type rows (20-27) with the same
# testing-type label already present in column E, add a brand new
# "Resource Usage Test" row's missing Expected Result / Test CaseType
# cells (D28 / E28), and update the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column B (Test Scenario Description) for rows 20-27 ---------------
# These rows previously had an empty column B; the same "testing type"
# text already shown in column E for that row gets duplicated into B.
$bValues = @{
    20 = "Exploratory Testing"
    21 = "Accessibility Testing"
    22 = "Localization Testing"
    23 = "API Testing"
    24 = "Database Testing"
    25 = "GUI Testing"
    26 = "Static Testing"
    27 = "Visual Testing"
}

foreach ($row in 20..27) {
    $cell = $ws.Range("B$row")
    $cell.Value = $bValues[$row]
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
}

# --- 2. New content for row 28 (Resource Usage Test) ----------------------
# Row 28 already had Test Case / Test Scenario Description / Test Case
# Description filled in; the Expected Result (D) and Test CaseType (E)
# cells are newly populated here.
$d28 = $ws.Range("D28")
$d28.Value = "There should be minimum utilization of resource."
$d28.Font.Name = "Arial"
$d28.Font.Size = 14
$d28.Borders.LineStyle = 1

$e28 = $ws.Range("E28")
$e28.Value = "Resource Usage Test"
$e28.Font.Name = "Arial"
$e28.Font.Size = 14
$e28.Borders.LineStyle = 1
$e28.HorizontalAlignment = -4108   # xlCenter
$e28.VerticalAlignment = -4160     # xlTop

# --- 3. Update the sheet's active selection --------------------------------
[void]$ws.Range("B20:B27").Select()
